$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.50080680847168
$ws.Range("B1").Value = 2.051681041717529
$ws.Range("C1").Value = 3.170416593551636
$ws.Range("D1").Value = 4.844185829162598
$ws.Range("E1").Value = 0.8570576310157776
